$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '71.427.24'
$ws.Range("E2").Value = '  +0.41%  '

$ws.Range("D3").Value = '2.570.26'
$ws.Range("E3").Value = '  +0.00%  '

$ws.Range("E4").Value = '  +0.07%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '583.58'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +0.24%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '172.64'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +0.94%  '

$ws.Range("E7").Value = '  +0.07%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.519'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  +1.32%  '

$ws.Range("B9").Value = 'LidoStakedEther'
$ws.Range("C9").Value = 'https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth'
$ws.Range("D9").Value = '2.568.50'
$ws.Range("E9").Value = '  -0.06%  '

$ws.Range("B10").Value = 'Dogecoin'
$ws.Range("C10").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.167'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -0.34%  '

$ws.Range("E11").Value = '  -0.41%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.361'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  +2.76%  '

$ws.Range("E13").Value = '  +1.68%  '

$ws.Range("D14").Value = '3.044.92'
$ws.Range("E14").Value = '  -0.33%  '

$ws.Range("D15").Value = '71.339.49'
$ws.Range("E15").Value = '  +0.61%  '

$ws.Range("E16").Value = '  -2.51%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '25.54'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  +1.01%  '

$ws.Range("D18").Value = '2.582.51'
$ws.Range("E18").Value = '  -1.53%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.67'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -1.92%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.97'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +3.59%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '358.62'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -1.71%  '

$ws.Range("E22").Value = '  -0.84%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '2.07'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +3.48%  '

$ws.Range("E24").Value = '  +0.03%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '70.72'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +0.34%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '4.12'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -0.88%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.18'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -0.86%  '

$ws.Range("D28").Value = '2.706.94'
$ws.Range("E28").Value = '  -1.75%  '

$ws.Range("E29").Value = '  +0.46%  '

$ws.Range("D30").Value = '0.0₃0928'
$ws.Range("E30").Value = '  -0.34%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '8.00'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  +2.40%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '477.37'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -2.02%  '

$ws.Range("E33").Value = '  -1.79%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.78'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +0.02%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.00'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +0.02%  '

$ws.Range("E36").Value = '  +3.48%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '157.52'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +0.83%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '18.91'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +0.29%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '19.13'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +1.41%  '

$ws.Range("E40").Value = '  +0.04%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '4.92'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +2.75%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.323'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +0.79%  '

$ws.Range("E43").Value = '  -3.99%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.39'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -3.62%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.18'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -11.34%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '38.78'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +0.37%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '146.50'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -0.62%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.542'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +1.81%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '3.57'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -0.38%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.63'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -0.50%  '

$ws.Range("E51").Value = '  +1.24%  '
